$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B106").Formula = "=_xlfn.STDEV.S(B2:B104)"
$ws.Range("C106").Formula = "=_xlfn.STDEV.S(C2:C104)"
$ws.Range("D106").Formula = "=_xlfn.STDEV.S(D2:D104)"
$ws.Range("E106").Formula = "=_xlfn.STDEV.S(E2:E104)"
$ws.Range("F106").Formula = "=_xlfn.STDEV.S(F2:F104)"

$ws.Range("B107").Formula = "=B106/SQRT(103)*1.96"
$ws.Range("C107").Formula = "=C106/SQRT(103)*1.96"
$ws.Range("D107").Formula = "=D106/SQRT(103)*1.96"
$ws.Range("E107").Formula = "=E106/SQRT(103)*1.96"
$ws.Range("F107").Formula = "=F106/SQRT(103)*1.96"

$ws.Range("I113").Select()
